$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.34"
$ws.Range("E2").Value = "'1.04%"
$ws.Range("E3").Value = "'0.25%"
$ws.Range("D4").Value = "'0.08026"
$ws.Range("E4").Value = "'-0.22%"
$ws.Range("D5").Value = "'1.985"
$ws.Range("E5").Value = "'5.18%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.329"
$ws.Range("E6").Value = "'0.74%"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").Value = "'2.613"
$ws.Range("E7").Value = "'-3.68%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9486"
$ws.Range("E8").Value = "'0.93%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1120"
$ws.Range("E9").Value = "'-4.38%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1854"
$ws.Range("E10").Value = "'-0.78%"
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").Value = "'10.65"
$ws.Range("E11").Value = "'24.91%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09831"
$ws.Range("E12").Value = "'-1.47%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04665"
$ws.Range("E13").Value = "'10.55%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1064"
$ws.Range("E14").Value = "'-0.16%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001260"
$ws.Range("E15").Value = "'-0.74%"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04059"
$ws.Range("E16").Value = "'-4.45%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005918"
$ws.Range("E17").Value = "'-0.12%"
$ws.Range("B18").Value = "OKB"
$ws.Range("C18").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D18").Value = "'43.84"
$ws.Range("E18").Value = "'-1.46%"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "'3.356"
$ws.Range("E19").Value = "'-6.52%"
$ws.Range("E20").Value = "'-0.26%"
$ws.Range("D21").Value = "'0.1405"
$ws.Range("E21").Value = "'3.66%"
$ws.Range("D22").Value = "'0.2542"
$ws.Range("E22").Value = "'-3.74%"
$ws.Range("D23").Value = "'0.001255"
$ws.Range("E23").Value = "'1.29%"
$ws.Range("D24").Value = "'0.004336"
$ws.Range("E24").Value = "'-2.73%"
$ws.Range("D25").Value = "'0.0001198"
$ws.Range("E25").Value = "'-0.40%"
$ws.Range("D26").Value = "'0.0003738"
$ws.Range("E26").Value = "'-6.53%"
$ws.Range("D38").Value = "'0.02584"
$ws.Range("E38").Value = "'-1.92%"
$ws.Range("D39").Value = "'0.05667"
$ws.Range("E39").Value = "'3.37%"
$ws.Range("D40").Value = "'0.007568"
$ws.Range("E40").Value = "'-1.61%"
$ws.Range("D41").Value = "'0.1395"
$ws.Range("E41").Value = "'0.07%"
$ws.Range("D42").Value = "'0.007535"
$ws.Range("E42").Value = "'3.74%"
$ws.Range("D43").Value = "'0.002012"
$ws.Range("E43").Value = "'-1.92%"
$ws.Range("D44").Value = "'0.008380"
$ws.Range("E44").Value = "'-3.57%"
$ws.Range("D45").Value = "'0.00007122"
$ws.Range("E45").Value = "'0.01%"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.41%"
$ws.Range("E47").Value = "'55.06%"
$ws.Range("D48").Value = "'0.003545"
$ws.Range("E48").Value = "'0.50%"
$ws.Range("D49").Value = "'0.00002097"
$ws.Range("E49").Value = "'-0.41%"
$ws.Range("D50").Value = "'0.0001997"
$ws.Range("E50").Value = "'-0.41%"
